# Update database to a new year (shift 1396..1400 reporting periods to 1397..1401)
# and refresh the underlying figures ("change read_price algorithm").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header rows: shift the five period labels forward by one year ---
$ws.Range("E8").Value  = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F8").Value  = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G8").Value  = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H8").Value  = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I8").Value  = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E24").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F24").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G24").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H24").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I24").Value = "دوازده ماهه منتهی به 1401/12"

# --- Data rows: new figures for columns E (1397) through I (1401) ---
$rowValues = @{
    10 = @(0,0,0,0,357543)
    11 = @(0,0,0,0,0)
    12 = @(0,0,0,11006,3041)
    13 = @(0,0,0,0,2988)
    14 = @(1155,3332,5451,12375,0)
    15 = @(300,1072,831,2307,373)
    16 = @(900,2113,2264,2608,10600)
    17 = @(9612,41853,57928,75139,109939)
    18 = @(0,0,0,0,0)
    19 = @(5555,20466,26537,39885,59848)
    20 = @(17522,68836,93011,143320,544332)
    26 = @(42,46,35,23,27)
    27 = @(150,162,197,204,212)
}

foreach ($r in $rowValues.Keys) {
    $rowNum = [int]$r
    $vals = $rowValues[$r]
    $col = 5
    foreach ($v in $vals) {
        $ws.Cells.Item($rowNum, $col).Value = $v
        $col = $col + 1
    }
}
